# Apply cryptos.xlsx price/volume refresh (GitHub Actions scheduled update)
# Prices with a single decimal separator are forced to stay as text
# (leading apostrophe) so Excel does not silently convert them to numbers,
# matching the source data which stores all Price/Volume cells as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.194.19"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "2.375.09"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'311.96"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'108.56"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "'40.97"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.0918"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "'0.978"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "2.740.46"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'15.29"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "2.388.90"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "45.274.76"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").Value = "'14.70"
$ws.Range("E19").Value = "  +12.34%  "
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'73.36"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'259.80"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'7.23"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  +9.36%  "
$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "'37.63"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'169.27"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = "  +3.14%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "'0.0354"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("D42").Value = "'99.78"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.892.41"
$ws.Range("E43").Value = "  +13.30%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'69.87"
$ws.Range("E44").Value = "  -2.64%  "
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "'12.86"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'81.56"
$ws.Range("E48").Value = "  +5.50%  "
$ws.Range("D49").Value = "'5.67"
$ws.Range("E49").Value = "  +7.42%  "
$ws.Range("D50").Value = "'112.52"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "'9.22"
$ws.Range("E51").Value = "  +2.66%  "
